# Update the "Metadata" sheet values: Version, Title, Date
$wb = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "7.0.0"
$meta.Range("B5").Value = "Snapshot Age In Years"
$meta.Range("B8").Value = "2022-09-01T20:48:10+00:00"

# The "Elements" sheet's "Short" column for the root Extension row reuses the
# same text as the Title above ("Age In Years" -> "Snapshot Age In Years").
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("K2").Value = "Snapshot Age In Years"
